$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$cleanStyle = $ws.Range("D4").Style

$ws.Range("D2").Value = '''63.872.45'
$ws.Range("D2").Style = $cleanStyle
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").Value = '''2.625.10'
$ws.Range("D3").Style = $cleanStyle
$ws.Range("E3").Value = '  -1.11%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''596.54'
$ws.Range("D5").Style = $cleanStyle
$ws.Range("E5").Value = '  -1.39%  '
$ws.Range("D6").Value = '''150.19'
$ws.Range("D6").Style = $cleanStyle
$ws.Range("E6").Value = '  +1.62%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("E9").Value = '  +0.56%  '
$ws.Range("D10").Value = '''5.66'
$ws.Range("D10").Style = $cleanStyle
$ws.Range("E10").Value = '  +1.02%  '
$ws.Range("D11").Value = '''0.382'
$ws.Range("D11").Style = $cleanStyle
$ws.Range("E11").Value = '  +3.39%  '
$ws.Range("D12").Value = '''0.151'
$ws.Range("D12").Style = $cleanStyle
$ws.Range("E12").Value = '  -1.24%  '
$ws.Range("D13").Value = '''27.69'
$ws.Range("D13").Style = $cleanStyle
$ws.Range("E13").Value = '  +0.38%  '
$ws.Range("D14").Value = '''3.096.78'
$ws.Range("D14").Style = $cleanStyle
$ws.Range("E14").Value = '  -1.04%  '
$ws.Range("D15").Value = '''63.693.84'
$ws.Range("D15").Style = $cleanStyle
$ws.Range("E15").Value = '  +0.18%  '
$ws.Range("E16").Value = '  +2.09%  '
$ws.Range("D17").Value = '''2.635.59'
$ws.Range("D17").Style = $cleanStyle
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("D18").Value = '''12.35'
$ws.Range("D18").Style = $cleanStyle
$ws.Range("E18").Value = '  +7.48%  '
$ws.Range("D19").Value = '''4.63'
$ws.Range("D19").Style = $cleanStyle
$ws.Range("E19").Value = '  +2.17%  '
$ws.Range("D20").Value = '''349.74'
$ws.Range("D20").Style = $cleanStyle
$ws.Range("E20").Value = '  +2.10%  '
$ws.Range("D21").Value = '''6.91'
$ws.Range("D21").Style = $cleanStyle
$ws.Range("E21").Value = '  -1.26%  '
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("E23").Value = '  +2.09%  '
$ws.Range("D24").Value = '''66.31'
$ws.Range("D24").Style = $cleanStyle
$ws.Range("E24").Value = '  -0.82%  '
$ws.Range("E25").Value = '  +13.23%  '
$ws.Range("E26").Value = '  +1.61%  '
$ws.Range("E27").Value = '  -1.43%  '
$ws.Range("B28").Value = 'Bittensor'
$ws.Range("C28").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D28").Value = '''562.36'
$ws.Range("D28").Style = $cleanStyle
$ws.Range("E28").Value = '  +1.63%  '
$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").Value = '''8.19'
$ws.Range("D29").Style = $cleanStyle
$ws.Range("E29").Value = '  +4.42%  '
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").Value = '''0.165'
$ws.Range("D30").Style = $cleanStyle
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("E32").Value = '  -0.24%  '
$ws.Range("D33").Value = '''0.0₃0846'
$ws.Range("D33").Style = $cleanStyle
$ws.Range("E33").Value = '  +3.38%  '
$ws.Range("D34").Value = '''1.74'
$ws.Range("D34").Style = $cleanStyle
$ws.Range("E34").Value = '  -1.13%  '
$ws.Range("D35").Value = '''5.22'
$ws.Range("D35").Style = $cleanStyle
$ws.Range("E35").Value = '  +0.43%  '
$ws.Range("D36").Value = '''169.60'
$ws.Range("D36").Style = $cleanStyle
$ws.Range("E36").Value = '  +1.40%  '
$ws.Range("D37").Value = '''0.409'
$ws.Range("D37").Style = $cleanStyle
$ws.Range("E37").Value = '  +0.48%  '
$ws.Range("E38").Value = '  -0.04%  '
$ws.Range("E39").Value = '  +1.83%  '
$ws.Range("D40").Value = '''19.39'
$ws.Range("D40").Style = $cleanStyle
$ws.Range("E40").Value = '  +1.26%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").Value = '''170.71'
$ws.Range("D41").Style = $cleanStyle
$ws.Range("E41").Value = '  +0.89%  '
$ws.Range("B42").Value = 'USDe'
$ws.Range("C42").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D42").Value = '''0.999'
$ws.Range("D42").Style = $cleanStyle
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").Value = '''39.82'
$ws.Range("D43").Style = $cleanStyle
$ws.Range("E43").Value = '  -0.22%  '
$ws.Range("E44").Value = '  +3.82%  '
$ws.Range("D45").Value = '''0.0598'
$ws.Range("D45").Style = $cleanStyle
$ws.Range("E45").Value = '  +3.51%  '
$ws.Range("D46").Value = '''21.46'
$ws.Range("D46").Style = $cleanStyle
$ws.Range("E46").Value = '  -5.45%  '
$ws.Range("D47").Value = '''0.631'
$ws.Range("D47").Style = $cleanStyle
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("D48").Value = '''0.0248'
$ws.Range("D48").Style = $cleanStyle
$ws.Range("E48").Value = '  -0.36%  '
$ws.Range("D49").Value = '''1.98'
$ws.Range("D49").Style = $cleanStyle
$ws.Range("E49").Value = '  +5.36%  '
$ws.Range("D50").Value = '''0.0968'
$ws.Range("D50").Style = $cleanStyle
$ws.Range("E50").Value = '  +0.43%  '
$ws.Range("D51").Value = '''19.23'
$ws.Range("D51").Style = $cleanStyle
$ws.Range("E51").Value = '  +1.72%  '
